$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17, J17, L17, N17
$ws.Range("H17").Value = 1050.125
$ws.Range("J17").Value = 1061.8096
$ws.Range("L17").Value = 3185.4288
$ws.Range("N17").Value = -3521.4288

# Row 29: H29, J29, L29, N29
$ws.Range("H29").Value = 3839.7144
$ws.Range("J29").Value = 6499.75
$ws.Range("L29").Value = 19499.25
$ws.Range("N29").Value = -20061.25

# Row 86: H86, I86, K86, M86
$ws.Range("H86").Value = 5030.3076
$ws.Range("I86").Value = 4565.143
$ws.Range("K86").Value = 4565.143
$ws.Range("M86").Value = -3442.143

# Row 87: H87, J87, L87, N87
$ws.Range("H87").Value = 87999.5
$ws.Range("J87").Value = 87999.5
$ws.Range("L87").Value = 87999.5
$ws.Range("N87").Value = -90495.5

# Row 89: H89, I89, K89, M89
$ws.Range("H89").Value = 5030.3076
$ws.Range("I89").Value = 4565.143
$ws.Range("K89").Value = 22825.715
$ws.Range("M89").Value = -17209.715

# Row 90: H90, J90, L90, N90
$ws.Range("H90").Value = 87999.5
$ws.Range("J90").Value = 87999.5
$ws.Range("L90").Value = 263998.5
$ws.Range("N90").Value = -276478.5

# Row 112: H112, J112, L112, N112
$ws.Range("H112").Value = 1229.9412
$ws.Range("J112").Value = 1223.9333
$ws.Range("L112").Value = 3671.7999
$ws.Range("N112").Value = -5887.7999

# Row 129: H129, I129, K129, M129
$ws.Range("H129").Value = 4912.8
$ws.Range("I129").Value = 4891
$ws.Range("K129").Value = 14673
$ws.Range("M129").Value = -9673

# Row 135: H135, I135, J135, K135, L135, M135, N135
$ws.Range("H135").Value = 3432.6924
$ws.Range("I135").Value = 323.27274
$ws.Range("J135").Value = 20534.5
$ws.Range("K135").Value = 2909.45466
$ws.Range("L135").Value = 184810.5
$ws.Range("M135").Value = -374.4546599999999
$ws.Range("N135").Value = -189880.5

# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 3338.0667
$ws.Range("I137").Value = 2875.9443
$ws.Range("J137").Value = 4031.25
$ws.Range("K137").Value = 8627.832900000001
$ws.Range("L137").Value = 12093.75
$ws.Range("M137").Value = -6077.832900000001
$ws.Range("N137").Value = -17193.75

# Row 138: H138, J138, L138, N138
$ws.Range("H138").Value = 3295.3262
$ws.Range("J138").Value = 3612.7144
$ws.Range("L138").Value = 10838.1432
$ws.Range("N138").Value = -21118.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32, I32, K32, M32
$ws.Range("H32").Value = 1938.1096
$ws.Range("I32").Value = 1206.1818
$ws.Range("K32").Value = 1206.1818
$ws.Range("M32").Value = -919.1818000000001

# Row 61: H61, I61, K61, M61
$ws.Range("H61").Value = 2221.5806
$ws.Range("I61").Value = 1946.8
$ws.Range("K61").Value = 1946.8
$ws.Range("M61").Value = -1734.8

# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 2208.5715
$ws.Range("I132").Value = 1835.5
$ws.Range("K132").Value = 5506.5
$ws.Range("M132").Value = -2976.5

# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 2221.5806
$ws.Range("I136").Value = 1946.8
$ws.Range("K136").Value = 5840.4
$ws.Range("M136").Value = -3290.4

$ws = $wb.Worksheets.Item("BSM")
# Row 94: H94, I94, J94, K94, L94, M94, N94
$ws.Range("H94").Value = 1205.4348
$ws.Range("I94").Value = 1191.1333
$ws.Range("J94").Value = 1232.25
$ws.Range("K94").Value = 1191.1333
$ws.Range("L94").Value = 1232.25
$ws.Range("M94").Value = -740.1333
$ws.Range("N94").Value = -2134.25

# Row 99: H99, I99, K99, M99
$ws.Range("H99").Value = 29489.75
$ws.Range("I99").Value = 32513.555
$ws.Range("K99").Value = 32513.555
$ws.Range("M99").Value = -31015.555

# Row 105: H105, J105, L105, N105
$ws.Range("H105").Value = 1223.8
$ws.Range("J105").Value = 1165
$ws.Range("L105").Value = 1165
$ws.Range("N105").Value = -4659

# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 1641.78
$ws.Range("I134").Value = 1130.7733
$ws.Range("J134").Value = 3174.8
$ws.Range("K134").Value = 3392.3199
$ws.Range("L134").Value = 9524.400000000001
$ws.Range("M134").Value = -857.3199000000004
$ws.Range("N134").Value = -14594.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31, J31, L31, N31
$ws.Range("H31").Value = 9842.35
$ws.Range("J31").Value = 8649.200000000001
$ws.Range("L31").Value = 8649.200000000001
$ws.Range("N31").Value = -9239.200000000001

# Row 34: H34, J34, L34, N34
$ws.Range("H34").Value = 9842.35
$ws.Range("J34").Value = 8649.200000000001
$ws.Range("L34").Value = 8649.200000000001
$ws.Range("N34").Value = -9053.200000000001

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1423.381
$ws.Range("I132").Value = 1423.381
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4270.143
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1740.143
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 76: H76, I76, J76, K76, L76, M76, N76
$ws.Range("H76").Value = 3802.6667
$ws.Range("I76").Value = 1706.5
$ws.Range("J76").Value = 7995
$ws.Range("K76").Value = 5119.5
$ws.Range("L76").Value = 23985
$ws.Range("M76").Value = -4736.5
$ws.Range("N76").Value = -24751

# Row 79: H79, I79, J79, K79, L79, M79, N79
$ws.Range("H79").Value = 3802.6667
$ws.Range("I79").Value = 1706.5
$ws.Range("J79").Value = 7995
$ws.Range("K79").Value = 5119.5
$ws.Range("L79").Value = 23985
$ws.Range("M79").Value = -3793.5
$ws.Range("N79").Value = -26637

# Row 87: H87, I87, K87, M87
$ws.Range("H87").Value = 15506.4
$ws.Range("I87").Value = 9000
$ws.Range("K87").Value = 27000
$ws.Range("M87").Value = -25752

# Row 90: H90, I90, K90, M90
$ws.Range("H90").Value = 15506.4
$ws.Range("I90").Value = 9000
$ws.Range("K90").Value = 81000
$ws.Range("M90").Value = -74760

# Row 121: H121, I121, J121, K121, L121, M121, N121
$ws.Range("H121").Value = 22223232
$ws.Range("I121").Value = 50000268
$ws.Range("J121").Value = 1602.6
$ws.Range("K121").Value = 150000804
$ws.Range("L121").Value = 4807.799999999999
$ws.Range("M121").Value = -149999494
$ws.Range("N121").Value = -7427.799999999999

# Row 131: H131, J131, L131, N131
$ws.Range("H131").Value = 1224.4318
$ws.Range("J131").Value = 1267.25
$ws.Range("L131").Value = 3801.75
$ws.Range("N131").Value = -13881.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80, I80, J80, K80, L80, M80, N80
$ws.Range("H80").Value = 318065.3
$ws.Range("I80").Value = 386773.53
$ws.Range("J80").Value = 169197.5
$ws.Range("K80").Value = 386773.53
$ws.Range("L80").Value = 169197.5
$ws.Range("M80").Value = -385775.53
$ws.Range("N80").Value = -171193.5

# Row 83: H83, I83, J83, K83, L83, M83, N83
$ws.Range("H83").Value = 318065.3
$ws.Range("I83").Value = 386773.53
$ws.Range("J83").Value = 169197.5
$ws.Range("K83").Value = 1933867.65
$ws.Range("L83").Value = 845987.5
$ws.Range("M83").Value = -1928875.65
$ws.Range("N83").Value = -855971.5

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 29423036
$ws.Range("I132").Value = 38469788
$ws.Range("J132").Value = 21097.125
$ws.Range("K132").Value = 115409364
$ws.Range("L132").Value = 63291.375
$ws.Range("M132").Value = -115406834
$ws.Range("N132").Value = -68351.375

$ws = $wb.Worksheets.Item("LTW")
# Row 93: H93, I93, J93, K93, L93, M93, N93
$ws.Range("H93").Value = 2635.353
$ws.Range("I93").Value = 2407.6155
$ws.Range("J93").Value = 3375.5
$ws.Range("K93").Value = 2407.6155
$ws.Range("L93").Value = 3375.5
$ws.Range("M93").Value = -1159.6155
$ws.Range("N93").Value = -5871.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 3506.389
$ws.Range("I132").Value = 2470.6
$ws.Range("K132").Value = 7411.799999999999
$ws.Range("M132").Value = -4881.799999999999
